# "Task 5: Handle Error Cases and wrap up" -> "Task 6: Handle Error Cases and wrap up"
#
# The source document has two "Task 5: " headings; only the second one (the
# final task heading, immediately followed by the run "Handle Error Cases and
# wrap up") is renumbered to "Task 6: ". We locate that paragraph by its full
# text (robust against any position/index drift), then rewrite it so the
# leading "Task 5: " run is split into three runs - "Task ", "6", ": " - each
# carrying the same run formatting the original single run had, exactly as
# in the target revision. The trailing "Handle Error Cases and wrap up" run
# is left untouched.

$d = $word.ActiveDocument

$target = $null
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text
    if ($text -ne $null -and $text.TrimEnd([char]13) -eq "Task 5: Handle Error Cases and wrap up") {
        $target = $para
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Task 5: Handle Error Cases and wrap up' paragraph"
}

$newParaXml = '<w:p w14:paraId="629CDA03" w14:textId="7F89801D" w:rsidR="00431BC8" w:rsidRPr="005077D4" w:rsidRDefault="00431BC8" w:rsidP="00431BC8"><w:pPr><w:spacing w:before="200"/><w:outlineLvl w:val="1"/><w:rPr><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r w:rsidRPr="005077D4"><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Task </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>6</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/><w:b/><w:bCs/><w:color w:val="38761D"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Handle Error Cases and wrap up</w:t></w:r></w:p>'

$target.Range.InsertXML($newParaXml)
